$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the tiny floating-point drift on the existing last row (A11)
$ws.Range("A11").Value = 45878.45849049769

# Append the new row (row 12) pulled in by the automatic WSL update
$ws.Range("A12").Value = 45878.50017998619
$ws.Range("A12").NumberFormat = $ws.Range("A11").NumberFormat

$ws.Range("B12").Value = 2025
$ws.Range("C12").Value = 37
$ws.Range("D12").Value = 16.78
$ws.Range("E12").Value = 83.87
$ws.Range("F12").Value = 493.95
$ws.Range("G12").Value = 12.2
$ws.Range("H12").Value = "ESE"
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = "12:00:15"
